$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header info ---
$ws.Range("B7").Value = "FAD"
$ws.Range("C7").Value = "PR No.:  2020-03-0124"

# F7 ("Date:" value cell) is the non-top-left half of a merged cell
# (E7:F7). A direct .Value assignment on a merged cell's non-anchor
# member is a no-op (matches real Excel), so route the write through a
# scratch cell + Copy/PasteSpecial, which *does* update it.
$ws.Range("Z100").Formula = '="March 25, 2020"'
$ws.Range("Z100").Copy()
$ws.Range("F7").PasteSpecial()
$ws.Range("Z100").Clear()

# --- Line item 1 (row 11) ---
$ws.Range("A11").Value = "S280"
$ws.Range("B11").Value = "box"
$ws.Range("C11").Value = "Map Pin/Pin Assrtd Colored`n"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 65
$ws.Range("F11").Value = 65

# --- Line item 2 (row 12) ---
$ws.Range("A12").Value = "S280"
$ws.Range("B12").Value = "box"
$ws.Range("C12").Value = "Map Pin/Pin Assrtd Colored`n"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 65
$ws.Range("F12").Value = 65

# --- Purpose (B37) ---
# "121212 " looks numeric, so a plain .Value assignment would silently
# coerce it to the Number 121212 (dropping the trailing space and the
# text type Excel originally stored). Force it through as text the same
# scratch-cell + PasteSpecial way, using a formula (whose result is
# already typed as Text) so no NumberFormat/style change leaks onto B37.
$ws.Range("Z100").Formula = '="121212 "'
$ws.Range("Z100").Copy()
$ws.Range("B37").PasteSpecial()
$ws.Range("Z100").Clear()

# --- Signatories ---
$ws.Range("B43").Value = "DR. CARINA S. CRUZ"
$ws.Range("B44").Value = "FAD Chief"
